$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B42: was stored as an inline string "3", should be a real number 3
$ws.Range("B42").Value = 3

# Add new row 43 with annotation data
$ws.Range("A43").Value = "Sunsi Wu"

# B43 needs to remain text "2" (not auto-converted to a number like Excel
# normally would do for a numeric-looking string). Force text entry by
# temporarily switching the cell to a text number format, then restore the
# default "Normal" style once the text value has been committed.
$b43 = $ws.Range("B43")
$b43.NumberFormat = "@"
$b43.Value = "2"
$b43.Style = "Normal"

$ws.Range("C43").Value = "again wrong"
$ws.Range("D43").Value = "FBK"
$ws.Range("E43").Value = "MET"
$ws.Range("F43").Value = "295c014b-37cb-453e-93b8-ae293d0d968b"
$ws.Range("G43").Value = "BkiIkBJ0b_annotated.xlsx"
$ws.Range("H43").Value = 'The other part of the criticism that we use a "straw man" is again wrong because we do not intend to show pathology with Mirowski et al. paper, experiments or claims.'
